$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 36, shifting existing rows 36..114 down to 37..115
$ws.Rows("36").Insert()

# Populate the newly inserted row 36 with the new data record
$ws.Cells.Item(36, 1).Value = 10
$ws.Cells.Item(36, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(36, 3).Value = "La Araucanía"
$ws.Cells.Item(36, 4).Value = 44526
$ws.Cells.Item(36, 5).Value = 9
$ws.Cells.Item(36, 6).Value = "Fruta"
$ws.Cells.Item(36, 7).Value = 100107
$ws.Cells.Item(36, 8).Value = "Otros"
$ws.Cells.Item(36, 9).Value = 100107002
$ws.Cells.Item(36, 10).Value = "Chirimoya"
$ws.Cells.Item(36, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(36, 12).Value = "Primera"
$ws.Cells.Item(36, 13).Value = 65
$ws.Cells.Item(36, 14).Value = 3000
$ws.Cells.Item(36, 15).Value = 3000
$ws.Cells.Item(36, 16).Value = 3000
$ws.Cells.Item(36, 17).Value = '$/kilo (en caja de 15 kilos)'
$ws.Cells.Item(36, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(36, 19).Value = 3000
$ws.Cells.Item(36, 20).Value = 1
